$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
# Row 2
$ws.Range("H2").Value = 1942.3182
$ws.Range("I2").Value = 556.44446
$ws.Range("J2").Value = 2901.7693
$ws.Range("K2").Value = 556.44446
$ws.Range("L2").Value = 2901.7693
$ws.Range("M2").Value = -443.44446
$ws.Range("N2").Value = -3127.7693
# Row 6
$ws.Range("H6").Value = 12066.667
$ws.Range("I6").Value = 12066.667
$ws.Range("K6").Value = 36200.001
$ws.Range("M6").Value = -36088.001
# Row 8
$ws.Range("H8").Value = 150.33333
$ws.Range("I8").Value = 225.0
$ws.Range("K8").Value = 675.0
$ws.Range("M8").Value = -536.0
# Row 15
$ws.Range("H15").Value = 628.9231
$ws.Range("I15").Value = 628.9231
$ws.Range("K15").Value = 1886.7693
$ws.Range("M15").Value = -1717.7693
# Row 41
$ws.Range("H41").Value = 382.22223
$ws.Range("I41").Value = 367.5
$ws.Range("K41").Value = 367.5
$ws.Range("M41").Value = 72.5
# Row 53
$ws.Range("H53").Value = 228.3158
$ws.Range("J53").Value = 397.66666
$ws.Range("L53").Value = 397.66666
$ws.Range("N53").Value = -1671.66666
# Row 55
$ws.Range("H55").Value = 256.85715
$ws.Range("I55").Value = 291.33334
$ws.Range("J55").Value = 50.0
$ws.Range("K55").Value = 291.33334
$ws.Range("L55").Value = 50.0
$ws.Range("M55").Value = -77.33334000000002
$ws.Range("N55").Value = -478.0
# Row 80
$ws.Range("H80").Value = 9880.6
$ws.Range("I80").Value = 9752.0
$ws.Range("J80").Value = 9966.333
$ws.Range("K80").Value = 29256.0
$ws.Range("L80").Value = 29898.999
$ws.Range("M80").Value = -28258.0
$ws.Range("N80").Value = -31894.999
# Row 83
$ws.Range("H83").Value = 9880.6
$ws.Range("I83").Value = 9752.0
$ws.Range("J83").Value = 9966.333
$ws.Range("K83").Value = 87768.0
$ws.Range("L83").Value = 89696.997
$ws.Range("M83").Value = -82776.0
$ws.Range("N83").Value = -99680.997
# Row 107
$ws.Range("H107").Value = 785.0
$ws.Range("I107").Value = 785.0
$ws.Range("K107").Value = 785.0
$ws.Range("M107").Value = 1135.0
# Row 111
$ws.Range("H111").Value = 1476.4166
$ws.Range("I111").Value = 604.2857
$ws.Range("J111").Value = 2697.4
$ws.Range("K111").Value = 1812.8571
$ws.Range("L111").Value = 8092.200000000001
$ws.Range("M111").Value = 1254.1429
$ws.Range("N111").Value = -14226.2
# Row 132
$ws.Range("H132").Value = 845.7
$ws.Range("I132").Value = 824.6111
$ws.Range("J132").Value = 1035.5
$ws.Range("K132").Value = 2473.8333
$ws.Range("L132").Value = 3106.5
$ws.Range("M132").Value = 56.16670000000022
$ws.Range("N132").Value = -8166.5

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
# Row 135
$ws.Range("H135").Value = 197499.0
$ws.Range("J135").Value = 197499.0
$ws.Range("L135").Value = 197499.0
$ws.Range("N135").Value = -207639.0

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
# Row 80
$ws.Range("H80").Value = 398.3158
$ws.Range("I80").Value = 634.25
$ws.Range("J80").Value = 226.72728
$ws.Range("K80").Value = 634.25
$ws.Range("L80").Value = 226.72728
$ws.Range("M80").Value = 363.75
$ws.Range("N80").Value = -2222.72728
# Row 83
$ws.Range("H83").Value = 398.3158
$ws.Range("I83").Value = 634.25
$ws.Range("J83").Value = 226.72728
$ws.Range("K83").Value = 3171.25
$ws.Range("L83").Value = 1133.6364
$ws.Range("M83").Value = 1820.75
$ws.Range("N83").Value = -11117.6364
# Row 135
$ws.Range("H135").Value = 99995.0
$ws.Range("J135").Value = 99995.0
$ws.Range("L135").Value = 99995.0
$ws.Range("N135").Value = -110135.0

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
# Row 22
$ws.Range("H22").Value = 200002.0
$ws.Range("I22").Value = 0.0
$ws.Range("J22").Value = 200002.0
$ws.Range("K22").Value = 0.0
$ws.Range("L22").Value = 200002.0
$ws.Range("M22").Value = $null
$ws.Range("N22").Value = -200702.0
# Row 131
$ws.Range("H131").Value = 69995.5
$ws.Range("J131").Value = 69995.5
$ws.Range("L131").Value = 69995.5
$ws.Range("N131").Value = -80075.5

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
# Row 4
$ws.Range("H4").Value = 6666770.0
$ws.Range("I4").Value = 6666770.0
$ws.Range("K4").Value = 20000310.0
$ws.Range("M4").Value = -20000198.0
# Row 6
$ws.Range("H6").Value = 100.57143
$ws.Range("I6").Value = 117.083336
$ws.Range("J6").Value = 1.5
$ws.Range("K6").Value = 351.250008
$ws.Range("L6").Value = 4.5
$ws.Range("M6").Value = -238.250008
$ws.Range("N6").Value = -230.5
# Row 98
$ws.Range("H98").Value = 210.33333
$ws.Range("I98").Value = 192.4
$ws.Range("K98").Value = 577.2
$ws.Range("M98").Value = 920.8
# Row 134
$ws.Range("H134").Value = 9052.223
$ws.Range("I134").Value = 1655.0
$ws.Range("J134").Value = 18298.75
$ws.Range("K134").Value = 4965.0
$ws.Range("L134").Value = 54896.25
$ws.Range("M134").Value = 105.0
$ws.Range("N134").Value = -65036.25

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
# Row 38
$ws.Range("I38").Value = 10000.0
$ws.Range("J38").Value = 0.0
$ws.Range("K38").Value = 10000.0
$ws.Range("L38").Value = 0.0
$ws.Range("M38").Value = -9537.0
$ws.Range("N38").Value = $null
# Row 132
$ws.Range("H132").Value = 924.5
$ws.Range("I132").Value = 899.6667
$ws.Range("K132").Value = 2699.0001
$ws.Range("M132").Value = -169.0001000000002

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
# Row 7
$ws.Range("H7").Value = 8630.789
$ws.Range("I7").Value = 0.0
$ws.Range("K7").Value = 0.0
$ws.Range("M7").Value = $null
# Row 22
$ws.Range("H22").Value = 1795.1428
$ws.Range("J22").Value = 1795.1428
$ws.Range("L22").Value = 1795.1428
$ws.Range("N22").Value = -2385.1428
# Row 27
$ws.Range("H27").Value = 1795.1428
$ws.Range("J27").Value = 1795.1428
$ws.Range("L27").Value = 1795.1428
$ws.Range("N27").Value = -2009.1428
# Row 40
$ws.Range("H40").Value = 4311.0
$ws.Range("I40").Value = 2600.0
$ws.Range("K40").Value = 2600.0
$ws.Range("M40").Value = -2464.0
# Row 46
$ws.Range("H46").Value = 64571.625
$ws.Range("J46").Value = 2000.0
$ws.Range("L46").Value = 2000.0
$ws.Range("N46").Value = -2376.0
# Row 55
$ws.Range("H55").Value = 247.22223
$ws.Range("I55").Value = 184.375
$ws.Range("K55").Value = 184.375
$ws.Range("M55").Value = -11.375
# Row 92
$ws.Range("H92").Value = 0.0
$ws.Range("J92").Value = 0.0
$ws.Range("L92").Value = 0.0
$ws.Range("N92").Value = $null
# Row 93
$ws.Range("H93").Value = 900.0
$ws.Range("I93").Value = 800.0
$ws.Range("J93").Value = 950.0
$ws.Range("K93").Value = 800.0
$ws.Range("L93").Value = 950.0
$ws.Range("M93").Value = 448.0
$ws.Range("N93").Value = -3446.0
# Row 98
$ws.Range("H98").Value = 0.0
$ws.Range("J98").Value = 0.0
$ws.Range("L98").Value = 0.0
$ws.Range("N98").Value = $null
# Row 126
$ws.Range("H126").Value = 8630.789
$ws.Range("I126").Value = 0.0
$ws.Range("K126").Value = 0.0
$ws.Range("M126").Value = $null
# Row 132
$ws.Range("H132").Value = 3244.6365
$ws.Range("I132").Value = 2637.2
$ws.Range("K132").Value = 7911.599999999999
$ws.Range("M132").Value = -5381.599999999999
# Row 136
$ws.Range("H136").Value = 3008.9443
$ws.Range("I136").Value = 2428.9167
$ws.Range("K136").Value = 7286.750100000001
$ws.Range("M136").Value = -4736.750100000001

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
# Row 29
$ws.Range("H29").Value = 5500.0
$ws.Range("I29").Value = 5000.0
$ws.Range("K29").Value = 5000.0
$ws.Range("M29").Value = -4710.0
# Row 113
$ws.Range("H113").Value = 339.5
$ws.Range("I113").Value = 302.66666
$ws.Range("J113").Value = 450.0
$ws.Range("K113").Value = 907.9999799999999
$ws.Range("L113").Value = 1350.0
$ws.Range("M113").Value = 1262.00002
$ws.Range("N113").Value = -5690.0
# Row 119
$ws.Range("H119").Value = 84999.75
$ws.Range("J119").Value = 84999.75
$ws.Range("L119").Value = 84999.75
$ws.Range("N119").Value = -94675.75
# Row 126
$ws.Range("H126").Value = 5375.8125
$ws.Range("I126").Value = 4439.875
$ws.Range("K126").Value = 13319.625
$ws.Range("M126").Value = -10849.625
